# Auto-generated edit script: update Leve profit sheets with refreshed market data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6333.2856
$ws.Range("I40").Value = 4653.2666
$ws.Range("J40").Value = 10533.333
$ws.Range("K40").Value = 4653.2666
$ws.Range("L40").Value = 10533.333
$ws.Range("M40").Value = -4478.2666
$ws.Range("N40").Value = -10883.333

$ws.Range("H76").Value = 20007936
$ws.Range("I76").Value = 55569556
$ws.Range("K76").Value = 55569556
$ws.Range("M76").Value = -55569241

$ws.Range("H79").Value = 20007936
$ws.Range("I79").Value = 55569556
$ws.Range("K79").Value = 55569556
$ws.Range("M79").Value = -55568464

$ws.Range("H98").Value = 3883.8235
$ws.Range("I98").Value = 2467.2727
$ws.Range("J98").Value = 6480.8335
$ws.Range("K98").Value = 2467.2727
$ws.Range("L98").Value = 6480.8335
$ws.Range("M98").Value = -969.2727
$ws.Range("N98").Value = -9476.833500000001

$ws.Range("H112").Value = 1601.3784
$ws.Range("I112").Value = 731.1111
$ws.Range("J112").Value = 1881.1072
$ws.Range("K112").Value = 2193.3333
$ws.Range("L112").Value = 5643.321599999999
$ws.Range("M112").Value = -1085.3333
$ws.Range("N112").Value = -7859.321599999999

$ws.Range("H122").Value = 3883.8235
$ws.Range("I122").Value = 2467.2727
$ws.Range("J122").Value = 6480.8335
$ws.Range("K122").Value = 7401.8181
$ws.Range("L122").Value = 19442.5005
$ws.Range("M122").Value = -4951.8181
$ws.Range("N122").Value = -24342.5005

$ws.Range("H137").Value = 5025.561
$ws.Range("I137").Value = 5654.4375
$ws.Range("J137").Value = 2789.5557
$ws.Range("K137").Value = 16963.3125
$ws.Range("L137").Value = 8368.667099999999
$ws.Range("M137").Value = -14413.3125
$ws.Range("N137").Value = -13468.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4354.9185
$ws.Range("I32").Value = 4748.525
$ws.Range("J32").Value = 2605.5557
$ws.Range("K32").Value = 4748.525
$ws.Range("L32").Value = 2605.5557
$ws.Range("M32").Value = -4461.525
$ws.Range("N32").Value = -3179.5557

$ws.Range("H61").Value = 457681.62
$ws.Range("I61").Value = 358778.47
$ws.Range("J61").Value = 630762.1
$ws.Range("K61").Value = 358778.47
$ws.Range("L61").Value = 630762.1
$ws.Range("M61").Value = -358566.47
$ws.Range("N61").Value = -631186.1

$ws.Range("H74").Value = 230304.77
$ws.Range("I74").Value = 334006.7
$ws.Range("J74").Value = 57468.223
$ws.Range("K74").Value = 334006.7
$ws.Range("L74").Value = 57468.223
$ws.Range("M74").Value = -333132.7
$ws.Range("N74").Value = -59216.223

$ws.Range("H77").Value = 230304.77
$ws.Range("I77").Value = 334006.7
$ws.Range("J77").Value = 57468.223
$ws.Range("K77").Value = 1670033.5
$ws.Range("L77").Value = 287341.115
$ws.Range("M77").Value = -1665665.5
$ws.Range("N77").Value = -296077.115

$ws.Range("H132").Value = 2005.0217
$ws.Range("I132").Value = 1335.8422
$ws.Range("J132").Value = 2475.926
$ws.Range("K132").Value = 4007.5266
$ws.Range("L132").Value = 7427.778
$ws.Range("M132").Value = -1477.5266
$ws.Range("N132").Value = -12487.778

$ws.Range("H136").Value = 457681.62
$ws.Range("I136").Value = 358778.47
$ws.Range("J136").Value = 630762.1
$ws.Range("K136").Value = 1076335.41
$ws.Range("L136").Value = 1892286.3
$ws.Range("M136").Value = -1073785.41
$ws.Range("N136").Value = -1897386.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3475.9822
$ws.Range("I134").Value = 3661.7441
$ws.Range("J134").Value = 2861.5386
$ws.Range("K134").Value = 10985.2323
$ws.Range("L134").Value = 8584.6158
$ws.Range("M134").Value = -8450.2323
$ws.Range("N134").Value = -13654.6158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2575.0908
$ws.Range("I31").Value = 776.6061
$ws.Range("J31").Value = 4373.5757
$ws.Range("K31").Value = 776.6061
$ws.Range("L31").Value = 4373.5757
$ws.Range("M31").Value = -481.6061
$ws.Range("N31").Value = -4963.5757

$ws.Range("H34").Value = 2575.0908
$ws.Range("I34").Value = 776.6061
$ws.Range("J34").Value = 4373.5757
$ws.Range("K34").Value = 776.6061
$ws.Range("L34").Value = 4373.5757
$ws.Range("M34").Value = -574.6061
$ws.Range("N34").Value = -4777.5757

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H132").Value = 2118.756
$ws.Range("I132").Value = 1408.2903
$ws.Range("J132").Value = 4321.2
$ws.Range("K132").Value = 4224.8709
$ws.Range("L132").Value = 12963.6
$ws.Range("M132").Value = -1694.8709
$ws.Range("N132").Value = -18023.6

$ws.Range("H134").Value = 1903.881
$ws.Range("I134").Value = 1522.1
$ws.Range("J134").Value = 2858.3333
$ws.Range("K134").Value = 4566.299999999999
$ws.Range("L134").Value = 8574.999899999999
$ws.Range("M134").Value = -2031.299999999999
$ws.Range("N134").Value = -13644.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1212.0435
$ws.Range("I131").Value = 1946.6666
$ws.Range("J131").Value = 1142.0793
$ws.Range("K131").Value = 5839.9998
$ws.Range("L131").Value = 3426.2379
$ws.Range("M131").Value = -799.9997999999996
$ws.Range("N131").Value = -13506.2379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4941.2
$ws.Range("I21").Value = 4853
$ws.Range("K21").Value = 4853
$ws.Range("M21").Value = -4680

$ws.Range("H30").Value = 4941.2
$ws.Range("I30").Value = 4853
$ws.Range("K30").Value = 4853
$ws.Range("M30").Value = -4748

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H132").Value = 2910397.2
$ws.Range("I132").Value = 4169866.2
$ws.Range("J132").Value = 3930.4614
$ws.Range("K132").Value = 12509598.6
$ws.Range("L132").Value = 11791.3842
$ws.Range("M132").Value = -12507068.6
$ws.Range("N132").Value = -16851.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6115.0186
$ws.Range("I132").Value = 1881.3334
$ws.Range("J132").Value = 8809.182000000001
$ws.Range("K132").Value = 5644.0002
$ws.Range("L132").Value = 26427.546
$ws.Range("M132").Value = -3114.0002
$ws.Range("N132").Value = -31487.546

$ws.Range("H136").Value = 4420.2144
$ws.Range("I136").Value = 2732
$ws.Range("J136").Value = 5814.826
$ws.Range("K136").Value = 8196
$ws.Range("L136").Value = 17444.478
$ws.Range("M136").Value = -5646
$ws.Range("N136").Value = -22544.478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2200.3489
$ws.Range("I132").Value = 1709.25
$ws.Range("J132").Value = 2820.6843
$ws.Range("K132").Value = 5127.75
$ws.Range("L132").Value = 8462.052899999999
$ws.Range("M132").Value = -2597.75
$ws.Range("N132").Value = -13522.0529

$ws.Range("H136").Value = 14603408
$ws.Range("I136").Value = 21981940
$ws.Range("J136").Value = 461219.22
$ws.Range("K136").Value = 65945820
$ws.Range("L136").Value = 1383657.66
$ws.Range("M136").Value = -65943270
$ws.Range("N136").Value = -1388757.66

Write-Host "Updated Leve profit sheets with refreshed market data."